# Remove the "farhan12" and "susipujiastuti" rows from the student table,
# leaving only the header, "bagustejo" and "adesusilo" rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (farhan12) entirely - this shifts everything below it up
$ws.Rows.Item(2).Delete()

# After the shift, "susipujiastuti" is now in row 2 - delete it too
$ws.Rows.Item(2).Delete()
